# ROLDANAS DISMAY.xlsx - fix(gui) step 1 and 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 1: move date forward one day (45308 -> 45309, i.e. 2024-01-17 -> 2024-01-18)
$ws.Range("A1").Value = 45309

# Updated prices (column D) for rows 28-31
$ws.Range("D28").Value = 230.1
$ws.Range("D29").Value = 300
$ws.Range("D30").Value = 336
$ws.Range("D31").Value = 422
